$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 3.223369029078222
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 8.656069925401464

# Row 3
$ws.Range("B3").Value = 3.272327238179451
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 6.15379541431027

# Row 4
$ws.Range("B4").Value = 0.01253208636536152
$ws.Range("C4").Value = 0.04103571897497393
$ws.Range("D4").Value = 0.1496068669990043
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("G4").Value = 0.7365606309410384

# Row 5
$ws.Range("B5").Value = 3.272327238179451
$ws.Range("C5").Value = 1.626987699542094
$ws.Range("D5").Value = 0.7210945179870265
$ws.Range("E5").Value = 0.5333859586016987
$ws.Range("G5").Value = 6.15379541431027
